# Donjons-et-Barons : "Amélioration du texte sur la diversité"
#
# 1) "types différents" -> "couleurs différentes" (three bullet points).
# 2) The explanatory sentence before the image table is replaced by a
#    shorter one.
# 3) "Recrutement spécial de la tuile-1" section reworded (heading + body).
# 4) The header's SAVEDATE field is refreshed so it reflects the (pinned)
#    save clock, mirroring what Word does automatically on save.

$d = $word.ActiveDocument

# --- 1) "types différents" -> "couleurs différentes" (bold kept) ---------
# Replace:=2 (wdReplaceAll) replaces every match in the range in one go,
# covering the three bullet points ("2/3/4 points pour N ...").
$rng = $d.Content
$rng.Find.Execute("types différents", $true, $false, $false, $false, $false, $true, 1, $false, "couleurs différentes", 2)

# --- 2) Replace the long explanatory sentence -----------------------------
$rng = $d.Content
$rng.Find.Execute( `
    "Le nombre de types est comptés comme si les tuiles gagnées étaient groupées en paquet, un paquet pour chaque type de tuiles : le nombre de types gagnés correspond au nombre de paquets constitués. En image :", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "En image, à une permutation près des 4 couleurs :", 2)

# --- 3) "Recrutement spécial de la tuile-1" section -----------------------
$rng = $d.Content
$rng.Find.Execute( `
    "Recrutement spécial de la tuile-1", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Recrutement spécial de la tuile T1", 2)

$rng = $d.Content
$rng.Find.Execute( `
    "La tuile-1, en abrégé « T1 », permet de recruter plus de troupe", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "La tuile « T1 » (à 1 point de victoire) permet de recruter plus de troupe", 2)

# --- 4) Refresh the cached SAVEDATE field result shown in the header ------
# (Word recomputes this "[yyyy-MMdd-HHmm]" field every time the file is
# saved; scope the Find to the header story so the body text is untouched.)
foreach ($sec in $d.Sections) {
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            $hdr.Range.Find.Execute("2024-0728-1107", $true, $false, $false, $false, $false, $true, 1, $false, "2024-0728-2118", 2)
        }
    }
}
